# "signed off time sheets" - fill in the supervisor's name and her
# sign-off (initials + date) on the timesheet now that it has been
# reviewed/approved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name (merged G6:I6, next to the "Supervisor Name:" label in E6)
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor sign-off row (mirrors the Employee sign-off in row 25):
# initials in the signature-line cell, and the date signed.
$ws.Range("A27").Value = "P.S"
$ws.Range("D27").Value = 41682
$ws.Range("D27").NumberFormat = $ws.Range("D25").NumberFormat
